# Auto-generated edit script applying cell value changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.699.97'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.468.20'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.97'
$ws.Range('E5').Value = '  +1.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.13'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.91'
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0855'
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.849.79'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.48'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.473.32'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.649.98'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.87'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.22'
$ws.Range('E22').Value = '  -2.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.50'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.78'
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.24'
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.18'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.34'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0763'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.03'
$ws.Range('E36').Value = '  -3.02%  '
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('E39').Value = '  +1.61%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.00'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -8.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.000.65'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.64'
$ws.Range('E45').Value = '  -1.82%  '
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.49'
$ws.Range('E47').Value = '  +4.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.730.59'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.67'
$ws.Range('E49').Value = '  +3.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.06'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '66.87'
$ws.Range('E51').Value = '  -0.30%  '
